$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-26 Friday" "2024-04-27 Saturday"

Replace-Text "77×24=" "57×44="
Replace-Text "58×14=" "38×24="
Replace-Text "64×43=" "66×36="
Replace-Text "90×76=" "33×70="
Replace-Text "78×42=" "36×35="

Replace-Text "13×66=" "40×26="
Replace-Text "51×58=" "95×85="
Replace-Text "77×82=" "83×99="
Replace-Text "47×66=" "54×68="
Replace-Text "18×15=" "91×67="

Replace-Text "97×51=" "30×64="
Replace-Text "70×45=" "96×81="
Replace-Text "97×53=" "44×11="
Replace-Text "42×85=" "89×57="
Replace-Text "51×18=" "20×66="

Replace-Text "65×78=" "98×74="
Replace-Text "96×73=" "62×51="
Replace-Text "17×85=" "13×17="
Replace-Text "51×24=" "77×39="
Replace-Text "28×53=" "70×49="

Replace-Text "12×95=" "50×87="
Replace-Text "81×68=" "57×35="
Replace-Text "48×48=" "19×17="
Replace-Text "24×63=" "37×30="
Replace-Text "69×21=" "85×37="
